$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Se agrega plataformas en packs/pos: nueva columna "Plataforma" (L)
$ws.Range("L1").Value = "Plataforma"
$ws.Range("L2").Value = "Saga Falabella"
$ws.Range("L3").Value = "Linio"

# Matches the style already used on row 3 (e.g. H3/I3) for the new cell
$ws.Range("L3").Style = $ws.Range("I3").Style

$ws.Range("K9").Select()
